$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (0603 LED): update link text + add hyperlink to the new "active part" (Harvatek) page ---
$ws.Range("C11").Value = "https://www.digikey.com/en/products/detail/harvatek-corporation/B19Y1USD-20C000113U1930/16602930"
$ws.Hyperlinks.Add($ws.Range("C11"), "https://www.digikey.com/en/products/detail/harvatek-corporation/B19Y1USD-20C000113U1930/16602930") | Out-Null

# --- Row 5 (2N7002): re-point the mouser link text (now with trailing space) + add hyperlink ---
$ws.Range("C5").Value = "https://www.mouser.com/ProductDetail/Nexperia/2N7002NXBKR?qs=%252B6g0mu59x7J2ddJstTJGkQ%3D%3D "
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.mouser.com/ProductDetail/Nexperia/2N7002NXBKR?qs=%252B6g0mu59x7J2ddJstTJGkQ%3D%3D") | Out-Null

# --- Reflect the author's final cursor position (was on A6, moved to C5) ---
$ws.Range("C5").Select() | Out-Null
